$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.036.91"
$ws.Range("E2").Value = "  +2.51%  "
$ws.Range("D3").Value = "'2.465.35"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'487.15"
$ws.Range("E5").Value = "  +4.40%  "
$ws.Range("D6").Value = "'145.16"
$ws.Range("E6").Value = "  +9.87%  "
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("D8").Value = "'0.508"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").Value = "'2.472.35"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'5.82"
$ws.Range("E10").Value = "  +9.28%  "
$ws.Range("D11").Value = "'0.0969"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "'2.896.00"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'56.056.33"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").Value = "'21.08"
$ws.Range("E16").Value = "  +6.79%  "
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "'2.476.46"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "'4.51"
$ws.Range("E19").Value = "  +6.56%  "
$ws.Range("D20").Value = "'10.05"
$ws.Range("E20").Value = "  +4.95%  "
$ws.Range("D21").Value = "'316.94"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'5.79"
$ws.Range("E23").Value = "  +7.55%  "
$ws.Range("D24").Value = "'58.44"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").Value = "'0.411"
$ws.Range("E25").Value = "  +6.08%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "'2.575.79"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("D29").Value = "'7.69"
$ws.Range("E29").Value = "  +7.36%  "
$ws.Range("D30").Value = "'0.0₃0782"
$ws.Range("E30").Value = "  +8.42%  "
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'147.78"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "'18.23"
$ws.Range("E33").Value = "  +2.47%  "
$ws.Range("E34").Value = "  +4.65%  "
$ws.Range("D35").Value = "'5.16"
$ws.Range("E35").Value = "  +2.58%  "
$ws.Range("E36").Value = "  +8.03%  "
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("D38").Value = "'0.860"
$ws.Range("E38").Value = "  +7.38%  "
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("E40").Value = "  +7.74%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "'0.0552"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("D43").Value = "'0.604"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("E44").Value = "  +6.47%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'261.36"
$ws.Range("E45").Value = "  +11.53%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0921"
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'4.71"
$ws.Range("E47").Value = "  +12.01%  "
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").Value = "'17.54"
$ws.Range("E50").Value = "  +5.51%  "
$ws.Range("D51").Value = "'1.871.38"
$ws.Range("E51").Value = "  -3.68%  "
